$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table originally had 6 "line" rows (line1..line6) followed by 8 "extr" rows
# (extr1..extr8). Two new line rows (line7, line8) were inserted right after
# line6, pushing the extr rows down by two rows. extr5 and extr6 also had
# their in_service flag swapped.

# First, extend the formatting (bold, centered, bordered index style) used by
# column A down into the two brand-new rows (16 and 17) by copying the format
# from the last existing data row.
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Data rows 8-17 (name, from_bus, to_bus, in_service) in final order.
$rows = @(
  @(6,  "line7", 14, 11, $true),
  @(7,  "line8", 16, 9,  $true),
  @(8,  "extr1", 5,  12, $true),
  @(9,  "extr2", 5,  9,  $true),
  @(10, "extr3", 10, 11, $false),
  @(11, "extr4", 7,  8,  $true),
  @(12, "extr5", 9,  11, $false),
  @(13, "extr6", 7,  11, $true),
  @(14, "extr7", 5,  7,  $true),
  @(15, "extr8", 8,  5,  $false)
)

$r = 8
foreach ($item in $rows) {
    $ws.Cells.Item($r, 1).Value = $item[0]
    $ws.Cells.Item($r, 2).Value = $item[1]
    $ws.Cells.Item($r, 3).Value = $item[2]
    $ws.Cells.Item($r, 4).Value = $item[3]
    $ws.Cells.Item($r, 5).Value = $item[4]
    $r = $r + 1
}
